# Refresh the crypto price / 1h-volume-change table (columns D and E,
# rows 2-51) with the latest scraped values from the Oct 12 2023
# GitHub Actions run.
#
# Column D holds price strings that often *look* numeric
# (e.g. "0.810", "1.10", "21.74") but must stay stored as TEXT, exactly
# like the source data (it even mixes thousands-separator dotted values
# like "26.852.33" with plain decimals). Assigning a bare numeric-looking
# string to Range.Value lets Excel's COM layer silently coerce it to a
# Double (dropping trailing zeros / introducing binary-float noise, e.g.
# "0.810" -> 0.81, "1.10" -> 1.1). To avoid that we prefix those values
# with a leading apostrophe, exactly as a user typing into the UI would,
# which forces Excel to keep them as text.
#
# Column E holds the "  +/-X.XX%  " strings (padded with spaces) which
# never parse as numbers, so they can be assigned directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '''26.852.33'
$ws.Range("E2").Value = '  -1.08%  '

# Row 3
$ws.Range("D3").Value = '''1.562.30'
$ws.Range("E3").Value = '  -0.03%  '

# Row 4
$ws.Range("E4").Value = '  -0.12%  '

# Row 5
$ws.Range("D5").Value = '''205.82'
$ws.Range("E5").Value = '  -0.31%  '

# Row 6
$ws.Range("D6").Value = '''0.489'
$ws.Range("E6").Value = '  -1.06%  '

# Row 7
$ws.Range("E7").Value = '  -0.07%  '

# Row 8
$ws.Range("D8").Value = '''21.74'
$ws.Range("E8").Value = '  -1.74%  '

# Row 9
$ws.Range("E9").Value = '  -0.20%  '

# Row 10
$ws.Range("E10").Value = '  -1.34%  '

# Row 11
$ws.Range("D11").Value = '''0.0864'
$ws.Range("E11").Value = '  +0.19%  '

# Row 12
$ws.Range("D12").Value = '''1.784.44'
$ws.Range("E12").Value = '  -0.03%  '

# Row 13
$ws.Range("D13").Value = '''1.564.01'
$ws.Range("E13").Value = '  +0.25%  '

# Row 14
$ws.Range("D14").Value = '''3.73'

# Row 15
$ws.Range("E15").Value = '  -0.40%  '

# Row 16
$ws.Range("D16").Value = '''26.872.16'

# Row 17
$ws.Range("D17").Value = '''61.25'
$ws.Range("E17").Value = '  -2.88%  '

# Row 18
$ws.Range("D18").Value = '''214.33'
$ws.Range("E18").Value = '  +1.12%  '

# Row 19
$ws.Range("D19").Value = '''7.36'
$ws.Range("E19").Value = '  +1.95%  '

# Row 20
$ws.Range("D20").Value = '''0.0₃0680'
$ws.Range("E20").Value = '  -1.16%  '

# Row 21
$ws.Range("E21").Value = '  -0.19%  '

# Row 22
$ws.Range("E22").Value = '  +0.28%  '

# Row 23
$ws.Range("D23").Value = '''9.15'
$ws.Range("E23").Value = '  -2.54%  '

# Row 24
$ws.Range("E24").Value = '  +0.94%  '

# Row 25
$ws.Range("D25").Value = '''153.93'
$ws.Range("E25").Value = '  +1.21%  '

# Row 26
$ws.Range("D26").Value = '''6.73'
$ws.Range("E26").Value = '  +2.49%  '

# Row 27
$ws.Range("E27").Value = '  +0.40%  '

# Row 28
$ws.Range("E28").Value = '  -0.14%  '

# Row 29
$ws.Range("E29").Value = '  -1.12%  '

# Row 30
$ws.Range("E30").Value = '  -0.06%  '

# Row 31
$ws.Range("D31").Value = '''1.10'
$ws.Range("E31").Value = '  -3.28%  '

# Row 32
$ws.Range("E32").Value = '  -0.22%  '

# Row 33
$ws.Range("D33").Value = '''1.402.34'

# Row 34
$ws.Range("E34").Value = '  -0.48%  '

# Row 35
$ws.Range("E35").Value = '  -1.49%  '

# Row 36
$ws.Range("E36").Value = '  -0.29%  '

# Row 37
$ws.Range("D37").Value = '''0.918'
$ws.Range("E37").Value = '  -2.73%  '

# Row 38
$ws.Range("E38").Value = '  -0.42%  '

# Row 39
$ws.Range("D39").Value = '''0.526'
$ws.Range("E39").Value = '  +0.95%  '

# Row 40
$ws.Range("D40").Value = '''0.810'
$ws.Range("E40").Value = '  -0.67%  '

# Row 41
$ws.Range("E41").Value = '  -0.13%  '

# Row 42
$ws.Range("D42").Value = '''0.997'
$ws.Range("E42").Value = '  +0.39%  '

# Row 43
$ws.Range("D43").Value = '''5.38'
$ws.Range("E43").Value = '  +3.26%  '

# Row 44
$ws.Range("E44").Value = '  +0.24%  '

# Row 45
$ws.Range("D45").Value = '''1.76'
$ws.Range("E45").Value = '  -0.90%  '

# Row 46
$ws.Range("D46").Value = '''63.15'
$ws.Range("E46").Value = '  -0.54%  '

# Row 47
$ws.Range("D47").Value = '''1.698.30'
$ws.Range("E47").Value = '  +0.14%  '

# Row 48
$ws.Range("D48").Value = '''86.22'
$ws.Range("E48").Value = '  +0.90%  '

# Row 49
$ws.Range("D49").Value = '''0.0504'
$ws.Range("E49").Value = '  +2.74%  '

# Row 50
$ws.Range("D50").Value = '''0.0₇0977'
$ws.Range("E50").Value = '  -1.97%  '

# Row 51
$ws.Range("E51").Value = '  +0.48%  '
